$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$requisitos = @(
  "LOB1003 -  Cálculo I  (Requisito)`n",
  "LOB1004 -  Cálculo II  (Requisito)`n",
  "LOB1006 -  Cálculo IV  (Requisito)`n",
  "LOB1008 -  Ciência, Tecnologia e Sociedade  (Requisito)`n",
  "LOB1012 -  Estatística  (Requisito)`n",
  "LOB1018 -  Física I  (Requisito)`n",
  "LOB1019 -  Física II  (Requisito)`n",
  "LOB1021 -  Física IV  (Requisito)`n",
  "LOB1036 -  Geometria Analítica  (Requisito)`n",
  "LOB1037 -  Àlgebra Linear  (Requisito)`n",
  "LOB1038 -  Física Experimental I  (Requisito)`n",
  "LOB1039 -  Física Experimental III  (Requisito)`n",
  "LOB1041 -  Física Experimental II  (Requisito)`n",
  "LOB1042 -  Física Experimental IV  (Requisito)`n",
  "LOB1045 -  Leitura e Produção de Textos Acadêmicos  (Requisito)`n",
  "LOB1052 -  Cálculo III  (Requisito)`n",
  "LOB1053 -  Física III  (Requisito)`n",
  "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)`n",
  "LOM3081 -  Introdução à Mecânica dos Sólidos  (Requisito)`n",
  "LOM3204 -  Desenho Técnico e Projeto Assistido por Computador  (Requisito)`n",
  "LOM3205 -  Eletromagnetismo  (Requisito)`n",
  "LOM3212 -  Fenômenos de Transporte A  (Requisito)`n",
  "LOM3218 -  Introdução à Engenharia Física  (Requisito)`n",
  "LOM3236 -  Processos de Fabricação  (Requisito)`n",
  "LOM3240 -  Química Inorgânica Fundamental e Aplicada  (Requisito)`n",
  "LOM3241 -  Química de Materiais  (Requisito)`n",
  "LOM3253 -  Física Matemática  (Requisito)`n",
  "LOM3257 -  Mecânica Clássica  (Requisito)`n",
  "LOM3260 -  Computação Científica em Python  (Requisito)`n",
  "LOM3261 -  Métodos Numéricos e Aplicações  (Requisito)`n",
  "LOM3262 -  Circuitos Elétricos  (Requisito)`n",
  "LOQ4095 -  Química Geral Experimental  (Requisito)`n",
  "LOQ4100 -  Fundamentos de Química para Engenharia I (Requisito)`n",
)

$startRow = 23
for ($i = 0; $i -lt $requisitos.Length; $i++) {
  $row = $startRow + $i
  $ws.Cells.Item($row, 2).Value = $requisitos[$i]
  $ws.Cells.Item($row, 3).Value = $requisitos[$i]
}

# Apply formatting (style + row height) to the newly added rows (44-55)
$ws.Range("B43:C43").Copy()
for ($row = 44; $row -le 55; $row++) {
  $ws.Range("B" + $row + ":C" + $row).PasteSpecial(-4122)
  $ws.Rows.Item($row).RowHeight = 30
}

Write-Host "Updated requisitos rows 23-55"
